$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.01581778927454889;  C = 0.1769645570152142;  D = 0.06154514893109243;  E = 0.2480829476830127;  F = 0.2569240011290105 }
    3  = @{ B = 0.02195526675411769;  C = 0.2163586422281868;  D = 0.08385254728376079;  E = 0.2895730430888911;  F = 0.3005296263947633 }
    4  = @{ B = -0.0193087473139601;  C = 0.2049953786970008;  D = 0.1019947020594539;   E = 0.3193660940980647;  F = 0.3329567938172152 }
    5  = @{ B = -0.07226700651545802; C = 0.1675871387573901;  D = 0.06404507904842406;  E = 0.2530712924225584;  F = 0.2543714109319102 }
    6  = @{ B = -0.1025096706313666;  C = 0.165819033310199;   D = 0.06384594992382764;  E = 0.252677561179911;   F = 0.2434422526798966 }
    7  = @{ B = -0.1148643701330149;  C = 0.1835058106740494;  D = 0.06418979330833566;  E = 0.2533570470863908;  F = 0.2395213268272089 }
    8  = @{ B = -0.09270844093133584; C = 0.2111731646910383;  D = 0.07863255137221827;  E = 0.2804149628179963;  F = 0.2899055632835616 }
    9  = @{ B = -0.2246702250377167;  C = 0.2246702250377167;  D = 0.1088070602894433;   E = 0.3298591521990004;  F = 0.2957964256146743 }
    10 = @{ B = -0.09344167408882031; C = 0.09344167408882031; D = 0.008731346456521313; E = 0.09344167408882031 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
